# Generate Report for Handoff
# Adds a new handed-off file (fe9c4b10-981f-41cf-92f9-11aa7b0ec7d4.md) as a
# new row (row 9) to every worksheet / table in the workbook:
#   - "Overview" sheet / "Overview" table
#   - "zh-cn" sheet     / "zh-cn" table
#   - "de-de" sheet     / "de-de" table

$wb = $excel.ActiveWorkbook

$fileId   = "fe9c4b10-981f-41cf-92f9-11aa7b0ec7d4"
$mdName   = "$fileId.md"
$mdPath   = "e2e\$fileId.md"
$ghBase   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8b2f7e6a3d4c9510eacb6f2491d7c835a0e6f9b2/e2e/$mdName"

$zhXlf    = "$fileId.f4a29374a690ae5303637942d3bbfed87949e7a6.zh-cn.xlf"
$deXlf    = "$fileId.f4a29374a690ae5303637942d3bbfed87949e7a6.de-de.xlf"

$zhHoDate = "2016-08-24 08:45:01"
$deHoDate = "2016-08-24 08:45:18"
$hoDate   = "2016-08-24 08:45:18"
$epoch    = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A9").Value2 = $mdName
$wsOverview.Range("B9").Value2 = $mdPath
$wsOverview.Range("C9").Value2 = ".md"
$wsOverview.Range("E9").Value2 = "Ready for handoff"
$wsOverview.Range("F9").Value2 = "Ready for handoff"
$wsOverview.Range("G9").Value2 = $hoDate
$wsOverview.Range("G9").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B9"), $ghBase, [Type]::Missing, [Type]::Missing, $mdPath) | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A9").Value2 = $mdName
$wsZh.Range("B9").Value2 = ".md"
$wsZh.Range("C9").Value2 = "Ready for handoff"
$wsZh.Range("D9").Value2 = "e2e"
$wsZh.Range("E9").Value2 = "ht"
$wsZh.Range("F9").Value2 = "'False"
$wsZh.Range("G9").Value2 = $zhXlf
$wsZh.Range("H9").Value2 = $zhHoDate
$wsZh.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K9").Value2 = $epoch
$wsZh.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("M9").Value2 = "'True"
$wsZh.Range("O9").Value2 = "'False"

$wsZh.Hyperlinks.Add($wsZh.Range("A9"), $ghBase, [Type]::Missing, [Type]::Missing, $mdName) | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A9").Value2 = $mdName
$wsDe.Range("B9").Value2 = ".md"
$wsDe.Range("C9").Value2 = "Ready for handoff"
$wsDe.Range("D9").Value2 = "e2e"
$wsDe.Range("E9").Value2 = "ht"
$wsDe.Range("F9").Value2 = "'False"
$wsDe.Range("G9").Value2 = $deXlf
$wsDe.Range("H9").Value2 = $deHoDate
$wsDe.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K9").Value2 = $epoch
$wsDe.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("M9").Value2 = "'True"
$wsDe.Range("O9").Value2 = "'False"

$wsDe.Hyperlinks.Add($wsDe.Range("A9"), $ghBase, [Type]::Missing, [Type]::Missing, $mdName) | Out-Null

$wb.Save()
